# This script applies the targeted numeric value updates described by the
# authoritative diff to the active worksheet of the workbook. Each entry in
# $changes specifies a Row (R), 1-based Column index (C), and the new
# numeric Value (V) that the cell must take.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{R=2; C=6; V=1.75},
    @{R=2; C=7; V=1.88},
    @{R=2; C=8; V=4.3},
    @{R=2; C=9; V=5.9},
    @{R=2; C=10; V=3.7},
    @{R=2; C=11; V=500},
    @{R=2; C=13; V=1.06},
    @{R=2; C=14; V=3.85},
    @{R=2; C=15; V=1.29},
    @{R=2; C=16; V=1.99},
    @{R=2; C=17; V=1.86},
    @{R=2; C=18; V=1.38},
    @{R=2; C=19; V=3.15},
    @{R=2; C=20; V=1.8},
    @{R=2; C=21; V=2.06},
    @{R=2; C=22; V=1.21},
    @{R=2; C=23; V=2.12},
    @{R=2; C=29; V=970},
    @{R=3; C=6; V=3.45},
    @{R=3; C=7; V=5.1},
    @{R=3; C=8; V=1.94},
    @{R=3; C=9; V=2.6},
    @{R=3; C=10; V=2.92},
    @{R=3; C=11; V=5},
    @{R=3; C=13; V=1.08},
    @{R=3; C=14; V=2.4},
    @{R=3; C=15; V=1.08},
    @{R=3; C=16; V=1.52},
    @{R=3; C=17; V=2.08},
    @{R=3; C=19; V=3},
    @{R=3; C=20; V=1.05},
    @{R=3; C=21; V=1.05},
    @{R=3; C=22; V=1.63},
    @{R=3; C=23; V=1.24},
    @{R=4; C=6; V=11.5},
    @{R=4; C=9; V=1.27},
    @{R=4; C=10; V=7.8},
    @{R=4; C=11; V=9},
    @{R=4; C=12; V=1.2},
    @{R=4; C=16; V=3.05},
    @{R=4; C=17; V=1.4},
    @{R=4; C=18; V=1.89},
    @{R=4; C=19; V=2.04},
    @{R=4; C=20; V=1.9},
    @{R=4; C=21; V=2},
    @{R=4; C=22; V=4.7},
    @{R=4; C=24; V=970},
    @{R=4; C=26; V=9.6},
    @{R=4; C=27; V=10},
    @{R=4; C=28; V=970},
    @{R=4; C=31; V=13},
    @{R=4; C=33; V=970},
    @{R=4; C=34; V=970},
    @{R=4; C=40; V=210},
    @{R=5; C=6; V=2.74},
    @{R=5; C=7; V=2.96},
    @{R=5; C=8; V=3},
    @{R=5; C=9; V=3.25},
    @{R=5; C=10; V=3.05},
    @{R=5; C=11; V=3.1},
    @{R=5; C=13; V=1.11},
    @{R=5; C=14; V=3.1},
    @{R=5; C=15; V=1.39},
    @{R=5; C=16; V=1.69},
    @{R=5; C=17; V=2.28},
    @{R=5; C=18; V=1.25},
    @{R=5; C=20; V=1.83},
    @{R=5; C=21; V=1.98},
    @{R=5; C=22; V=1.45},
    @{R=5; C=23; V=1.51},
    @{R=5; C=24; V=11.5},
    @{R=5; C=25; V=13.5},
    @{R=5; C=26; V=20},
    @{R=5; C=28; V=10.5},
    @{R=5; C=29; V=7},
    @{R=5; C=33; V=13},
    @{R=5; C=34; V=970},
    @{R=5; C=36; V=48},
    @{R=5; C=39; V=150},
    @{R=5; C=40; V=40},
    @{R=5; C=41; V=46},
    @{R=6; C=6; V=1.38},
    @{R=6; C=7; V=1.44},
    @{R=6; C=8; V=8.6},
    @{R=6; C=11; V=5.8},
    @{R=6; C=12; V=1.27},
    @{R=6; C=14; V=4.8},
    @{R=6; C=15; V=1.21},
    @{R=6; C=16; V=2.3},
    @{R=6; C=17; V=1.62},
    @{R=6; C=19; V=2.52},
    @{R=6; C=20; V=2},
    @{R=6; C=21; V=1.92},
    @{R=6; C=22; V=1.11},
    @{R=6; C=23; V=3.2},
    @{R=6; C=24; V=24},
    @{R=6; C=25; V=38},
    @{R=6; C=26; V=95},
    @{R=6; C=27; V=330},
    @{R=6; C=28; V=10},
    @{R=6; C=29; V=13},
    @{R=6; C=30; V=36},
    @{R=6; C=31; V=150},
    @{R=6; C=32; V=9.6},
    @{R=6; C=33; V=11},
    @{R=6; C=34; V=27},
    @{R=6; C=35; V=120},
    @{R=6; C=36; V=12.5},
    @{R=6; C=37; V=15.5},
    @{R=6; C=38; V=36},
    @{R=6; C=39; V=140},
    @{R=6; C=41; V=170},
    @{R=7; C=6; V=5.4},
    @{R=7; C=7; V=5.7},
    @{R=7; C=8; V=1.78},
    @{R=7; C=9; V=1.79},
    @{R=7; C=12; V=1.44},
    @{R=7; C=14; V=3.45},
    @{R=7; C=15; V=1.37},
    @{R=7; C=16; V=1.83},
    @{R=7; C=17; V=2.1},
    @{R=7; C=18; V=1.32},
    @{R=7; C=19; V=3.9},
    @{R=7; C=21; V=1.93},
    @{R=7; C=22; V=2.24},
    @{R=7; C=23; V=1.21},
    @{R=7; C=27; V=18.5},
    @{R=7; C=28; V=17},
    @{R=7; C=31; V=970},
    @{R=7; C=33; V=22},
    @{R=7; C=35; V=42},
    @{R=7; C=37; V=90},
    @{R=7; C=41; V=13},
    @{R=8; C=6; V=1.78},
    @{R=8; C=20; V=1.61},
    @{R=8; C=22; V=1.23},
    @{R=8; C=24; V=970},
    @{R=8; C=25; V=970},
    @{R=8; C=30; V=970},
    @{R=9; C=6; V=4.8},
    @{R=9; C=8; V=1.79},
    @{R=9; C=9; V=1.87},
    @{R=9; C=10; V=3.65},
    @{R=9; C=11; V=3.95},
    @{R=9; C=21; V=1.96},
    @{R=9; C=22; V=2.1},
    @{R=10; C=8; V=2.26},
    @{R=10; C=11; V=3.45},
    @{R=10; C=12; V=1.48},
    @{R=10; C=24; V=13},
    @{R=10; C=29; V=8.4},
    @{R=10; C=34; V=21},
    @{R=11; C=8; V=1.64},
    @{R=11; C=9; V=1.77},
    @{R=11; C=11; V=4.9},
    @{R=11; C=16; V=2.44},
    @{R=11; C=17; V=1.49},
    @{R=11; C=23; V=1.19},
    @{R=11; C=25; V=970},
    @{R=12; C=8; V=1.4},
    @{R=12; C=9; V=1.42},
    @{R=12; C=11; V=5.4},
    @{R=12; C=15; V=1.3},
    @{R=12; C=16; V=2.02},
    @{R=12; C=17; V=1.88},
    @{R=12; C=33; V=38},
    @{R=13; C=7; V=4.9},
    @{R=13; C=9; V=1.81},
    @{R=13; C=11; V=5.1},
    @{R=13; C=12; V=1.21},
    @{R=13; C=16; V=2.92},
    @{R=13; C=18; V=1.78},
    @{R=13; C=20; V=1.52},
    @{R=13; C=22; V=2.22},
    @{R=13; C=26; V=1000},
    @{R=13; C=28; V=32},
    @{R=13; C=29; V=12.5},
    @{R=13; C=30; V=11.5},
    @{R=13; C=31; V=1000},
    @{R=13; C=33; V=21},
    @{R=13; C=34; V=18},
    @{R=13; C=35; V=28},
    @{R=13; C=38; V=46},
    @{R=13; C=40; V=32},
    @{R=13; C=41; V=6.4},
    @{R=14; C=17; V=1.76},
    @{R=14; C=19; V=2.92},
    @{R=14; C=20; V=1.68},
    @{R=14; C=21; V=2.32},
    @{R=14; C=24; V=20},
    @{R=14; C=28; V=12.5},
    @{R=14; C=29; V=9.6},
    @{R=14; C=30; V=16},
    @{R=14; C=31; V=40},
    @{R=14; C=32; V=16.5},
    @{R=14; C=33; V=12},
    @{R=14; C=34; V=18},
    @{R=14; C=36; V=29},
    @{R=14; C=40; V=15},
    @{R=15; C=12; V=1.21},
    @{R=15; C=14; V=6.4},
    @{R=15; C=15; V=1.15},
    @{R=15; C=18; V=1.74},
    @{R=15; C=19; V=2.18},
    @{R=15; C=37; V=48},
    @{R=15; C=38; V=46},
    @{R=16; C=6; V=2.3},
    @{R=16; C=7; V=2.44},
    @{R=16; C=8; V=2.76},
    @{R=16; C=9; V=3.05},
    @{R=16; C=10; V=4.1},
    @{R=16; C=12; V=1.21},
    @{R=16; C=14; V=6.2},
    @{R=16; C=16; V=2.96},
    @{R=16; C=17; V=1.43},
    @{R=16; C=18; V=1.78},
    @{R=16; C=19; V=2.04},
    @{R=16; C=20; V=1.51},
    @{R=16; C=21; V=2.84},
    @{R=16; C=22; V=1.5},
    @{R=16; C=23; V=1.69},
    @{R=16; C=40; V=11.5},
    @{R=17; C=6; V=3.7},
    @{R=17; C=8; V=1.81},
    @{R=17; C=9; V=1.89},
    @{R=17; C=10; V=4.7},
    @{R=17; C=14; V=8},
    @{R=17; C=16; V=3.35},
    @{R=17; C=19; V=1.86},
    @{R=17; C=20; V=1.45},
    @{R=17; C=21; V=2.92},
    @{R=17; C=22; V=2.12},
    @{R=17; C=24; V=970},
    @{R=17; C=25; V=970},
    @{R=17; C=28; V=970},
    @{R=18; C=6; V=3.1},
    @{R=18; C=7; V=3.25},
    @{R=18; C=8; V=2.58},
    @{R=18; C=9; V=2.66},
    @{R=18; C=11; V=3.3},
    @{R=18; C=12; V=1.49},
    @{R=18; C=14; V=3.3},
    @{R=18; C=16; V=1.76},
    @{R=18; C=17; V=2.24},
    @{R=18; C=18; V=1.28},
    @{R=18; C=19; V=4.2},
    @{R=18; C=20; V=1.91},
    @{R=18; C=21; V=2.02},
    @{R=18; C=22; V=1.6},
    @{R=18; C=24; V=12},
    @{R=18; C=26; V=15.5},
    @{R=18; C=30; V=12},
    @{R=18; C=31; V=32},
    @{R=18; C=32; V=20},
    @{R=18; C=33; V=14},
    @{R=18; C=34; V=18.5},
    @{R=18; C=36; V=55},
    @{R=18; C=37; V=40},
    @{R=18; C=39; V=120},
    @{R=18; C=40; V=44},
    @{R=18; C=41; V=28},
    @{R=19; C=20; V=2.18},
    @{R=20; C=15; V=1.5},
    @{R=21; C=6; V=4.5},
    @{R=21; C=7; V=4.6},
    @{R=21; C=10; V=3.45},
    @{R=21; C=12; V=1.51},
    @{R=21; C=16; V=1.71},
    @{R=21; C=21; V=1.9},
    @{R=21; C=22; V=1.95},
    @{R=21; C=23; V=1.27},
    @{R=21; C=24; V=10.5},
    @{R=21; C=27; V=23},
    @{R=21; C=31; V=23},
    @{R=21; C=32; V=30},
    @{R=21; C=35; V=46},
    @{R=21; C=37; V=65},
    @{R=22; C=6; V=1.83},
    @{R=22; C=7; V=1.84},
    @{R=22; C=10; V=3.6},
    @{R=22; C=11; V=3.65},
    @{R=22; C=19; V=4.5},
    @{R=22; C=20; V=2.16},
    @{R=22; C=21; V=1.82},
    @{R=22; C=23; V=2.18},
    @{R=22; C=25; V=15},
    @{R=22; C=29; V=8},
    @{R=22; C=32; V=9.199999999999999},
    @{R=22; C=37; V=21},
    @{R=22; C=40; V=16},
    @{R=22; C=41; V=140},
    @{R=23; C=6; V=2.6},
    @{R=23; C=9; V=3.25},
    @{R=23; C=13; V=1.11},
    @{R=23; C=14; V=3.1},
    @{R=23; C=16; V=1.69},
    @{R=23; C=17; V=2.4},
    @{R=23; C=18; V=1.25},
    @{R=23; C=24; V=10},
    @{R=23; C=27; V=55},
    @{R=23; C=28; V=8.800000000000001},
    @{R=23; C=30; V=13.5},
    @{R=23; C=41; V=48},
    @{R=24; C=6; V=2.2},
    @{R=24; C=9; V=4.7},
    @{R=24; C=10; V=2.94},
    @{R=24; C=15; V=1.64},
    @{R=24; C=21; V=1.67},
    @{R=24; C=28; V=6.6},
    @{R=24; C=32; V=12.5},
    @{R=24; C=40; V=38},
    @{R=25; C=7; V=15.5},
    @{R=25; C=10; V=6.2},
    @{R=25; C=12; V=1.36},
    @{R=25; C=15; V=1.28},
    @{R=25; C=17; V=1.84},
    @{R=25; C=18; V=1.42},
    @{R=25; C=19; V=3.1},
    @{R=25; C=23; V=1.06},
    @{R=25; C=25; V=7.6},
    @{R=25; C=28; V=36},
    @{R=25; C=29; V=14.5},
    @{R=25; C=33; V=55},
    @{R=25; C=37; V=360},
    @{R=25; C=38; V=270},
    @{R=25; C=39; V=330},
    @{R=25; C=40; V=700},
    @{R=25; C=41; V=5.8},
    @{R=28; C=6; V=3.6},
    @{R=28; C=7; V=3.7},
    @{R=28; C=8; V=2.54},
    @{R=28; C=9; V=2.6},
    @{R=28; C=11; V=2.98},
    @{R=28; C=14; V=2.24},
    @{R=28; C=16; V=1.4},
    @{R=28; C=19; V=7.6},
    @{R=28; C=20; V=2.5},
    @{R=28; C=23; V=1.37},
    @{R=28; C=25; V=6.8},
    @{R=28; C=27; V=42},
    @{R=28; C=30; V=14.5},
    @{R=28; C=37; V=75},
    @{R=28; C=39; V=370},
    @{R=28; C=41; V=55},
    @{R=29; C=7; V=3.65},
    @{R=29; C=8; V=2.6},
    @{R=29; C=10; V=2.86},
    @{R=29; C=12; V=1.68},
    @{R=29; C=19; V=6.6},
    @{R=29; C=23; V=1.38},
    @{R=31; C=17; V=2.56},
    @{R=33; C=6; V=2.08},
    @{R=33; C=7; V=2.14},
    @{R=33; C=8; V=4.8},
    @{R=33; C=23; V=1.87},
    @{R=33; C=30; V=22},
    @{R=33; C=31; V=100},
    @{R=33; C=40; V=36},
    @{R=34; C=11; V=3.65},
    @{R=34; C=41; V=290},
    @{R=35; C=14; V=2.9},
    @{R=35; C=17; V=2.4},
    @{R=35; C=19; V=4.6},
    @{R=35; C=27; V=29}
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg.R, $chg.C).Value = $chg.V
}
